$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells that changed value only ---
$ws.Range("F2").Value = 45069
$ws.Range("F6").Value = 45065

# --- Row 7: data previously shown for "WHAIR" is replaced by new ticker "TFFIF" data ---
$ws.Range("A7").Value = "TFFIF"
$ws.Range("B7").Value = 0.1032
$ws.Range("C7").Value = 20000
$ws.Range("D7").Value = 2064
$ws.Range("E7").Value = 1857.6
$ws.Range("F7").Value = 45077
$ws.Range("G7").Value = 45096
$ws.Range("H7").Value = 153000
$ws.Range("I7").Value = 1.21
$ws.Range("J7").Value = 1

# --- Row 8: brand-new row holding "WHAIR" data, inherit date format from row above ---
$ws.Range("F8").NumberFormat = $ws.Range("F7").NumberFormat
$ws.Range("G8").NumberFormat = $ws.Range("G7").NumberFormat

$ws.Range("A8").Value = "WHAIR"
$ws.Range("B8").Value = 0.1369
$ws.Range("C8").Value = 50000
$ws.Range("D8").Value = 6845
$ws.Range("E8").Value = 6160.5
$ws.Range("F8").Value = 45068
$ws.Range("G8").Value = 45099
$ws.Range("H8").Value = 435000
$ws.Range("I8").Value = 1.42
$ws.Range("J8").Value = 1

# --- Row 9: brand-new row holding "WHART" data (old row 8 content, shifted down,
#     with xdate/paiddate/actual refreshed), inherit date format from row above ---
$ws.Range("F9").NumberFormat = $ws.Range("F7").NumberFormat
$ws.Range("G9").NumberFormat = $ws.Range("G7").NumberFormat

$ws.Range("A9").Value = "WHART"
$ws.Range("B9").Value = 0.1915
$ws.Range("C9").Value = 30000
$ws.Range("D9").Value = 5745
$ws.Range("E9").Value = 5170.5
$ws.Range("F9").Value = 45069
$ws.Range("G9").Value = 45086
$ws.Range("H9").Value = 351000
$ws.Range("I9").Value = 1.47
$ws.Range("J9").Value = 1
